$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per repulled data
$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -6
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = 7
